$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'knee pads snowboard'
$ws.Cells.Item(2, 1).Value = 'youth under armour basketball tights'
$ws.Cells.Item(3, 1).Value = 'black basketball rim'
$ws.Cells.Item(4, 1).Value = 'mcdavid youth knee pads basketball'
$ws.Cells.Item(5, 1).Value = 'nike compression pants men'
$ws.Cells.Item(6, 1).Value = 'asics knee pad'
$ws.Cells.Item(7, 1).Value = 'nike leggings basketball'
$ws.Cells.Item(8, 1).Value = 'compression tights with knee pads'
$ws.Cells.Item(9, 1).Value = 'combat knee pads'
$ws.Cells.Item(10, 1).Value = 'uflex knee compression'
$ws.Cells.Item(11, 1).Value = 'damascus knee pads'
$ws.Cells.Item(12, 1).Value = 'mcdavid youth basketball knee pads'
$ws.Cells.Item(13, 1).Value = 'mens pants black'
$ws.Cells.Item(14, 1).Value = 'basketball pants girls'
$ws.Cells.Item(15, 1).Value = 'compression capri leggings'
$ws.Cells.Item(16, 1).Value = 'padded tights'
$ws.Cells.Item(17, 1).Value = 'padded tights men'
$ws.Cells.Item(18, 1).Value = 'crye knee pads'
$ws.Cells.Item(19, 1).Value = 'dancing knee pads'
$ws.Cells.Item(20, 1).Value = 'dodgeball knee pads'
$ws.Cells.Item(21, 1).Value = 'uflex athletics knee compression sleeve'
$ws.Cells.Item(22, 1).Value = 'gform knee pads'
$ws.Cells.Item(23, 1).Value = 'nike tights for men'
$ws.Cells.Item(24, 1).Value = 'insertable knee pads'
$ws.Cells.Item(25, 1).Value = 'knee pads for scootering'
$ws.Cells.Item(26, 1).Value = 'youth compression pants with pads'
$ws.Cells.Item(27, 1).Value = 'red knee pads'
$ws.Cells.Item(28, 1).Value = 'black volleyball knee pads'
$ws.Cells.Item(29, 1).Value = 'men basketball pants'
$ws.Cells.Item(30, 1).Value = 'compression tights with pads'
$ws.Cells.Item(31, 1).Value = 'compression leggings with knee pads'
$ws.Cells.Item(32, 1).Value = 'ski knee pads'
$ws.Cells.Item(33, 1).Value = 'dead on knee pads'
$ws.Cells.Item(34, 1).Value = 'black pads'
$ws.Cells.Item(35, 1).Value = 'padded basketball compression pants'
$ws.Cells.Item(36, 1).Value = 'mens padded leggings'
$ws.Cells.Item(37, 1).Value = 'knee pads for skating'
$ws.Cells.Item(38, 1).Value = 'padded compression tights'
$ws.Cells.Item(39, 1).Value = 'basketball knee pads pants'
$ws.Cells.Item(40, 1).Value = 'compression pants womens'
$ws.Cells.Item(41, 1).Value = 'military knee pad'
$ws.Cells.Item(42, 1).Value = 'drskin mens compression pants'
$ws.Cells.Item(43, 1).Value = 'padded compression pants men basketball'
$ws.Cells.Item(44, 1).Value = 'snowboard knee pad'
$ws.Cells.Item(45, 1).Value = 'double knee pads'
$ws.Cells.Item(46, 1).Value = 'knee pads light'
$ws.Cells.Item(47, 1).Value = 'knee pad protection'
$ws.Cells.Item(48, 1).Value = 'leggings with knee pads women'
$ws.Cells.Item(49, 1).Value = 'padded tights men basketball'
$ws.Cells.Item(50, 1).Value = 'mens leggings with knee pads'
$ws.Cells.Item(51, 1).Value = 'knee pad under pants'
$ws.Cells.Item(52, 1).Value = 'football pants adult xxl'
$ws.Cells.Item(53, 1).Value = 'the best leggings'
$ws.Cells.Item(54, 1).Value = 'sliding shorts mens'
$ws.Cells.Item(55, 1).Value = 'snowboarding padded shorts men'
$ws.Cells.Item(56, 1).Value = 'flexible knee pads'
$ws.Cells.Item(57, 1).Value = 'basketball shorts for men pack'
$ws.Cells.Item(58, 1).Value = 'leggings capri'
$ws.Cells.Item(59, 1).Value = 'leggings spandex'
$ws.Cells.Item(60, 1).Value = 'floor knee pads'
$ws.Cells.Item(61, 1).Value = 'gym pads'
$ws.Cells.Item(62, 1).Value = 'pants with knee pads kids'
$ws.Cells.Item(63, 1).Value = 'basketball pants men'
$ws.Cells.Item(64, 1).Value = 'volleyball knee pads nike black'
$ws.Cells.Item(65, 1).Value = 'tights men'
$ws.Cells.Item(66, 1).Value = 'hayabusa compression pants'
$ws.Cells.Item(67, 1).Value = 'tough knee pads'
$ws.Cells.Item(68, 1).Value = 'knee pad leggings'
$ws.Cells.Item(69, 1).Value = 'goalkeeper pant'
$ws.Cells.Item(70, 1).Value = 'knee compression men'
$ws.Cells.Item(71, 1).Value = 'wrestling knee sleeves'
$ws.Cells.Item(72, 1).Value = 'padded knee sleeves for basketball'
$ws.Cells.Item(73, 1).Value = 'knee compression sleeve for squats'
$ws.Cells.Item(74, 1).Value = 'youth leg sleeves for basketball'
$ws.Cells.Item(75, 1).Value = 'tactical knee pads'
$ws.Cells.Item(76, 1).Value = 'knee pads paintball'
$ws.Cells.Item(77, 1).Value = 'bike knee pads'
$ws.Cells.Item(78, 1).Value = 'knee pads biking'
$ws.Cells.Item(79, 1).Value = 'goalie knee pads'
$ws.Cells.Item(80, 1).Value = 'knee pads compression sleeve'
$ws.Cells.Item(81, 1).Value = 'pant with knee pads'
$ws.Cells.Item(82, 1).Value = 'long knee pads'
$ws.Cells.Item(83, 1).Value = 'knee pads for girls'
$ws.Cells.Item(84, 1).Value = 'mens pants with knee pads'
$ws.Cells.Item(85, 1).Value = 'maroon knee pads'
$ws.Cells.Item(86, 1).Value = 'pantalones con rodilleras'
$ws.Cells.Item(87, 1).Value = 'pantalon con rodilleras'
$ws.Cells.Item(88, 1).Value = 'knee padded pants men'
$ws.Cells.Item(89, 1).Value = 'knee pad pants men'
$ws.Cells.Item(90, 1).Value = 'baseball sliding pants mens'
$ws.Cells.Item(91, 1).Value = 'men basketball knee pads'
$ws.Cells.Item(92, 1).Value = 'baseball pants mens knee'
$ws.Cells.Item(93, 1).Value = 'compression leggings basketball'
$ws.Cells.Item(94, 1).Value = 'black mens baseball pants'
$ws.Cells.Item(95, 1).Value = 'basketball pants for men'
$ws.Cells.Item(96, 1).Value = 'work pants with knee pads'
$ws.Cells.Item(97, 1).Value = 'knee pad for basketball youth'
$ws.Cells.Item(98, 1).Value = 'basketball tights for men'
$ws.Cells.Item(99, 1).Value = 'black athletic pants men'
$ws.Cells.Item(100, 1).Value = 'youth basketball clothes'
